{"js": "// The edit updates the date shown in the title paragraph and replaces\n// every arithmetic expression in the 20x5 answers table with a new one,\n// cell by cell (by position), while leaving all formatting untouched.\n\nconst body = context.document.body;\n\n// --- 1. Title paragraph (e.g. \"2024-10-23 Wednesday\" -> \"2024-10-24 Thursday\") ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titleParagraph = paragraphs.items[0];\nconst newTitle = \"2024-10-24 Thursday\";\n// Replacing the text of the existing range keeps the run's formatting\n// (font, size, etc.) intact instead of inserting a differently formatted run.\ntitleParagraph.getRange().insertText(newTitle, Word.InsertLocation.replace);\n\n// --- 2. The answers table: overwrite all 100 cells, row-major order ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// New value for every cell, outer array = rows (top to bottom), inner\n// array = columns (left to right), matching the table's existing 20x5 shape.\nconst newValues = [\n  [\"71-38=33\", \"20+70=90\", \"97-15=82\", \"85+5=90\", \"87-23=64\"],\n  [\"30+5=35\", \"59-48=11\", \"70-40=30\", \"0+70=70\", \"57+1=58\"],\n  [\"78+20=98\", \"81-13=68\", \"28+17=45\", \"91+4=95\", \"14+26=40\"],\n  [\"18+14=32\", \"34+32=66\", \"17+45=62\", \"60-57=3\", \"52+45=97\"],\n  [\"56-37=19\", \"99-82=17\", \"59-4=55\", \"19+14=33\", \"25+41=66\"],\n  [\"68-30=38\", \"50-50=0\", \"2+84=86\", \"72+4=76\", \"31+56=87\"],\n  [\"58+23=81\", \"28+6=34\", \"78-69=9\", \"12+74=86\", \"22+56=78\"],\n  [\"16+6=22\", \"43-42=1\", \"85-27=58\", \"97-43=54\", \"61-35=26\"],\n  [\"62-50=12\", \"14+23=37\", \"63-6=57\", \"83-41=42\", \"60-39=21\"],\n  [\"48+40=88\", \"1+9=10\", \"40-0=40\", \"51+12=63\", \"86-49=37\"],\n  [\"60-1=59\", \"10+78=88\", \"77+20=97\", \"91-35=56\", \"65+3=68\"],\n  [\"44-2=42\", \"40-14=26\", \"66-48=18\", \"56+20=76\", \"72-71=1\"],\n  [\"58+36=94\", \"35+43=78\", \"13+74=87\", \"46-15=31\", \"42+46=88\"],\n  [\"31+56=87\", \"88+9=97\", \"30+44=74\", \"24-3=21\", \"34+6=40\"],\n  [\"63-56=7\", \"36+41=77\", \"31+0=31\", \"72+8=80\", \"75-49=26\"],\n  [\"86-56=30\", \"44+12=56\", \"15+67=82\", \"8+13=21\", \"61+31=92\"],\n  [\"85+5=90\", \"4+10=14\", \"38+21=59\", \"85-31=54\", \"85-83=2\"],\n  [\"2+23=25\", \"52+15=67\", \"99-87=12\", \"81-11=70\", \"15+31=46\"],\n  [\"80+12=92\", \"54+8=62\", \"53+46=99\", \"87-27=60\", \"22+59=81\"],\n  [\"35-10=25\", \"47+38=85\", \"78-4=74\", \"10-9=1\", \"39+8=47\"]\n];\n\n// Assigning .values rewrites each cell's text in place and keeps the\n// original cell/paragraph/run formatting (fonts, size, alignment, etc.).\ntable.values = newValues;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date title (first paragraph in the document).\n$d.Paragraphs.Item(1).Range.Text = '2024-10-24 Thursday'\n\n# Update each cell of the 20x5 answers table, row by row, column by\n# column, to match the new values exactly (positional replacement,\n# since several cells share identical original text).\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @('71-38=33', '20+70=90', '97-15=82', '85+5=90', '87-23=64'),\n    @('30+5=35', '59-48=11', '70-40=30', '0+70=70', '57+1=58'),\n    @('78+20=98', '81-13=68', '28+17=45', '91+4=95', '14+26=40'),\n    @('18+14=32', '34+32=66', '17+45=62', '60-57=3', '52+45=97'),\n    @('56-37=19', '99-82=17', '59-4=55', '19+14=33', '25+41=66'),\n    @('68-30=38', '50-50=0', '2+84=86', '72+4=76', '31+56=87'),\n    @('58+23=81', '28+6=34', '78-69=9', '12+74=86', '22+56=78'),\n    @('16+6=22', '43-42=1', '85-27=58', '97-43=54', '61-35=26'),\n    @('62-50=12', '14+23=37', '63-6=57', '83-41=42', '60-39=21'),\n    @('48+40=88', '1+9=10', '40-0=40', '51+12=63', '86-49=37'),\n    @('60-1=59', '10+78=88', '77+20=97', '91-35=56', '65+3=68'),\n    @('44-2=42', '40-14=26', '66-48=18', '56+20=76', '72-71=1'),\n    @('58+36=94', '35+43=78', '13+74=87', '46-15=31', '42+46=88'),\n    @('31+56=87', '88+9=97', '30+44=74', '24-3=21', '34+6=40'),\n    @('63-56=7', '36+41=77', '31+0=31', '72+8=80', '75-49=26'),\n    @('86-56=30', '44+12=56', '15+67=82', '8+13=21', '61+31=92'),\n    @('85+5=90', '4+10=14', '38+21=59', '85-31=54', '85-83=2'),\n    @('2+23=25', '52+15=67', '99-87=12', '81-11=70', '15+31=46'),\n    @('80+12=92', '54+8=62', '53+46=99', '87-27=60', '22+59=81'),\n    @('35-10=25', '47+38=85', '78-4=74', '10-9=1', '39+8=47')\n)\n\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $row = $newValues[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n\n"}
